$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5872273333333333
$ws.Range("H2").Value = 1.761682
$ws.Range("I2").Value = 0.07214749117712801
$ws.Range("J2").Value = 0.07214749117712801
$ws.Range("M2").Value = 0.902915
$ws.Range("N2").Value = 2.708745
$ws.Range("O2").Value = 0.1151652421792931
$ws.Range("P2").Value = 0.1151652421792931
$ws.Range("Q2").Value = 0.5302163676766667
$ws.Range("R2").Value = 4.77194730909
$ws.Range("S2").Value = 0.008308883294042358
$ws.Range("T2").Value = 0.008308883294042358
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5872273333333333
$ws.Range("H3").Value = 1.761682
$ws.Range("I3").Value = 0.07214749117712801
$ws.Range("J3").Value = 0.07214749117712801
$ws.Range("O3").Value = 0.5443877317615758
$ws.Range("P3").Value = 0.5443877317615758
$ws.Range("Q3").Value = 2.506340283581333
$ws.Range("R3").Value = 22.557062552232
$ws.Range("S3").Value = 0.03927620907420502
$ws.Range("T3").Value = 0.03927620907420502
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5872273333333333
$ws.Range("H4").Value = 1.761682
$ws.Range("I4").Value = 0.07214749117712801
$ws.Range("J4").Value = 0.07214749117712801
$ws.Range("M4").Value = 2.440679
$ws.Range("N4").Value = 7.322037
$ws.Range("O4").Value = 0.311304373187858
$ws.Range("P4").Value = 0.311304373187858
$ws.Range("Q4").Value = 1.433233420692666
$ws.Range("R4").Value = 12.899100786234
$ws.Range("S4").Value = 0.02245982951797235
$ws.Range("T4").Value = 0.02245982951797235
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5872273333333333
$ws.Range("H5").Value = 1.761682
$ws.Range("I5").Value = 0.07214749117712801
$ws.Range("J5").Value = 0.07214749117712801
$ws.Range("M5").Value = 0.2284833333333333
$ws.Range("N5").Value = 0.68545
$ws.Range("O5").Value = 0.02914265287127302
$ws.Range("P5").Value = 0.02914265287127302
$ws.Range("Q5").Value = 0.1341716585444445
$ws.Range("R5").Value = 1.2075449269
$ws.Range("S5").Value = 0.002102569290908274
$ws.Range("T5").Value = 0.002102569290908275
$ws.Range("I6").Value = 0.927852508822872
$ws.Range("J6").Value = 0.927852508822872
$ws.Range("M6").Value = 0.902915
$ws.Range("N6").Value = 2.708745
$ws.Range("O6").Value = 0.1151652421792931
$ws.Range("P6").Value = 0.1151652421792931
$ws.Range("Q6").Value = 6.818845381053333
$ws.Range("R6").Value = 61.36960842948
$ws.Range("S6").Value = 0.1068563588852507
$ws.Range("T6").Value = 0.1068563588852507
$ws.Range("I7").Value = 0.927852508822872
$ws.Range("J7").Value = 0.927852508822872
$ws.Range("O7").Value = 0.5443877317615758
$ws.Range("P7").Value = 0.5443877317615758
$ws.Range("S7").Value = 0.5051115226873708
$ws.Range("T7").Value = 0.5051115226873708
$ws.Range("I8").Value = 0.927852508822872
$ws.Range("J8").Value = 0.927852508822872
$ws.Range("M8").Value = 2.440679
$ws.Range("N8").Value = 7.322037
$ws.Range("O8").Value = 0.311304373187858
$ws.Range("P8").Value = 0.311304373187858
$ws.Range("Q8").Value = 18.43209241820533
$ws.Range("R8").Value = 165.888831763848
$ws.Range("S8").Value = 0.2888445436698856
$ws.Range("T8").Value = 0.2888445436698857
$ws.Range("I9").Value = 0.927852508822872
$ws.Range("J9").Value = 0.927852508822872
$ws.Range("M9").Value = 0.2284833333333333
$ws.Range("N9").Value = 0.68545
$ws.Range("O9").Value = 0.02914265287127302
$ws.Range("P9").Value = 0.02914265287127302
$ws.Range("Q9").Value = 1.725514054088889
$ws.Range("R9").Value = 15.5296264868
$ws.Range("S9").Value = 0.02704008358036474
$ws.Range("T9").Value = 0.02704008358036475
